$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 1.75
$ws.Range("O2").Value = 1.06
$ws.Range("Q2").Value = 1.97
$ws.Range("R2").Value = 1.23
$ws.Range("S2").Value = 1.05
$ws.Range("F3").Value = 5.6
$ws.Range("H3").Value = 1.57
$ws.Range("K3").Value = 5.6
$ws.Range("L3").Value = 1.35
$ws.Range("S4").Value = 2.66
$ws.Range("V4").Value = 1.3
$ws.Range("AD4").Value = 16
$ws.Range("S5").Value = 1.8
$ws.Range("T5").Value = 2.46
$ws.Range("AA5").Value = 8.199999999999999
$ws.Range("AL5").Value = 450
$ws.Range("Q6").Value = 1.62
$ws.Range("R6").Value = 1.53
$ws.Range("S6").Value = 2.34
$ws.Range("AM7").Value = 85
$ws.Range("V8").Value = 1.58
$ws.Range("AM8").Value = 1000
$ws.Range("R9").Value = 1.63
$ws.Range("T9").Value = 1.64
$ws.Range("W12").Value = 1.69
$ws.Range("F13").Value = 1.65
$ws.Range("G13").Value = 1.73
$ws.Range("J13").Value = 4.5
$ws.Range("P13").Value = 2.44
$ws.Range("R13").Value = 1.58
$ws.Range("W13").Value = 2.36
$ws.Range("AM13").Value = 75
$ws.Range("F15").Value = 1.94
$ws.Range("I15").Value = 4.3
$ws.Range("U15").Value = 2.28
$ws.Range("AD15").Value = 19.5
$ws.Range("N16").Value = 5.1
$ws.Range("L17").Value = 1.27
$ws.Range("AH17").Value = 17
$ws.Range("L18").Value = 1.56
$ws.Range("Q18").Value = 2.58
$ws.Range("Z18").Value = 980
$ws.Range("J19").Value = 3.3
$ws.Range("S19").Value = 1.98
$ws.Range("G20").Value = 3.45
$ws.Range("K20").Value = 3.25
$ws.Range("L20").Value = 1.53
$ws.Range("W20").Value = 1.41
$ws.Range("X20").Value = 9.199999999999999
$ws.Range("AG20").Value = 15
$ws.Range("G21").Value = 2.28
$ws.Range("L21").Value = 1.45
$ws.Range("N21").Value = 3.45
$ws.Range("P21").Value = 1.85
$ws.Range("Q21").Value = 2.02
$ws.Range("W21").Value = 1.78
$ws.Range("P22").Value = 1.66
$ws.Range("I23").Value = 4.2
$ws.Range("V23").Value = 1.32
$ws.Range("AK23").Value = 32
$ws.Range("F24").Value = 1.89
$ws.Range("G24").Value = 1.93
$ws.Range("N24").Value = 3.3
$ws.Range("O24").Value = 1.39
$ws.Range("Q24").Value = 2.14
$ws.Range("S24").Value = 3.95
$ws.Range("W24").Value = 2.06
$ws.Range("H25").Value = 2.48
$ws.Range("F29").Value = 1.95
$ws.Range("G29").Value = 2.34
$ws.Range("H29").Value = 3.3
$ws.Range("J29").Value = 3.5
$ws.Range("K29").Value = 4.4
$ws.Range("M29").Value = 1.05
$ws.Range("N29").Value = 3.65
$ws.Range("O29").Value = 1.3
$ws.Range("P29").Value = 1.8
$ws.Range("S29").Value = 2.92
$ws.Range("W29").Value = 1.75
$ws.Range("Y29").Value = 21
$ws.Range("AC29").Value = 12.5
$ws.Range("AD29").Value = 23
$ws.Range("AF29").Value = 18.5
$ws.Range("H31").Value = 1.69
$ws.Range("I31").Value = 1.7
$ws.Range("F32").Value = 1.34
$ws.Range("G32").Value = 1.35
$ws.Range("H32").Value = 10
$ws.Range("I32").Value = 10.5
$ws.Range("T32").Value = 1.93
$ws.Range("U32").Value = 2
$ws.Range("W32").Value = 3.85
$ws.Range("X32").Value = 24
$ws.Range("P34").Value = 3.45
$ws.Range("F35").Value = 1.58
$ws.Range("G35").Value = 1.59
$ws.Range("H35").Value = 6.2
$ws.Range("I35").Value = 6.4
$ws.Range("K35").Value = 4.8
$ws.Range("N35").Value = 5.6
$ws.Range("P35").Value = 2.56
$ws.Range("Q35").Value = 1.61
$ws.Range("U35").Value = 2.34
$ws.Range("AA35").Value = 190
$ws.Range("AD35").Value = 23
$ws.Range("AM35").Value = 85
$ws.Range("Z36").Value = 16
$ws.Range("F37").Value = 1.7
$ws.Range("G37").Value = 1.71
$ws.Range("H37").Value = 5.5
$ws.Range("I37").Value = 5.6
$ws.Range("R37").Value = 1.56
$ws.Range("S37").Value = 2.68
$ws.Range("V37").Value = 1.21
$ws.Range("W37").Value = 2.4
$ws.Range("AM37").Value = 85
$ws.Range("AN37").Value = 7.8
$ws.Range("F39").Value = 1.69
$ws.Range("H39").Value = 4.5
$ws.Range("I39").Value = 5
$ws.Range("N41").Value = 2.94
$ws.Range("P41").Value = 1.63
$ws.Range("R41").Value = 1.24
$ws.Range("T41").Value = 2.34
$ws.Range("V41").Value = 1.15
$ws.Range("AM41").Value = 260
$ws.Range("G42").Value = 2.4
